$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "29.382.93"
Set-TextCell $ws.Range("E2") "  +0.50%  "

Set-TextCell $ws.Range("D3") "1.873.86"
Set-TextCell $ws.Range("E3") "  +0.78%  "

Set-TextCell $ws.Range("D4") "0.9993"

Set-TextCell $ws.Range("D5") "0.7138"
Set-TextCell $ws.Range("E5") "  +1.75%  "

Set-TextCell $ws.Range("D6") "238.76"
Set-TextCell $ws.Range("E6") "  +0.51%  "

Set-TextCell $ws.Range("D7") "0.9994"
Set-TextCell $ws.Range("E7") "  -0.13%  "

Set-TextCell $ws.Range("D8") "0.07816"
Set-TextCell $ws.Range("E8") "  -5.66%  "

Set-TextCell $ws.Range("D9") "0.3070"
Set-TextCell $ws.Range("E9") "  +1.18%  "

Set-TextCell $ws.Range("D10") "25.38"
Set-TextCell $ws.Range("E10") "  +9.23%  "

Set-TextCell $ws.Range("D11") "0.08199"
Set-TextCell $ws.Range("E11") "  +0.23%  "

Set-TextCell $ws.Range("D12") "1.862.63"
Set-TextCell $ws.Range("E12") "  -0.18%  "

Set-TextCell $ws.Range("D13") "5.250"
Set-TextCell $ws.Range("E13") "  +1.52%  "

Set-TextCell $ws.Range("D14") "0.7227"
Set-TextCell $ws.Range("E14") "  +1.49%  "

Set-TextCell $ws.Range("D15") "89.44"
Set-TextCell $ws.Range("E15") "  +0.39%  "

Set-TextCell $ws.Range("D16") "29.511.70"
Set-TextCell $ws.Range("E16") "  +0.86%  "

Set-TextCell $ws.Range("D17") "5.821"
Set-TextCell $ws.Range("E17") "  +0.77%  "

Set-TextCell $ws.Range("D18") "243.10"
Set-TextCell $ws.Range("E18") "  +2.56%  "

Set-TextCell $ws.Range("D19") "0.000007848"
Set-TextCell $ws.Range("E19") "  -0.03%  "

Set-TextCell $ws.Range("D20") "13.32"
Set-TextCell $ws.Range("E20") "  -0.18%  "

Set-TextCell $ws.Range("B21") "WrappedliquidstakedEther2.0"
Set-TextCell $ws.Range("C21") "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextCell $ws.Range("D21") "2.149.33"
Set-TextCell $ws.Range("E21") "  +1.59%  "

Set-TextCell $ws.Range("B22") "Dai"
Set-TextCell $ws.Range("C22") "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws.Range("D22") "0.9999"
Set-TextCell $ws.Range("E22") "  -0.09%  "

Set-TextCell $ws.Range("D23") "0.9996"
Set-TextCell $ws.Range("E23") "  -0.16%  "

Set-TextCell $ws.Range("D24") "7.746"
Set-TextCell $ws.Range("E24") "  +4.13%  "

Set-TextCell $ws.Range("B25") "Stellar"
Set-TextCell $ws.Range("C25") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws.Range("D25") "0.1470"
Set-TextCell $ws.Range("E25") "  +1.89%  "

Set-TextCell $ws.Range("B26") "Monero"
Set-TextCell $ws.Range("C26") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D26") "162.14"
Set-TextCell $ws.Range("E26") "  +0.22%  "

Set-TextCell $ws.Range("D27") "8.974"
Set-TextCell $ws.Range("E27") "  +0.06%  "

Set-TextCell $ws.Range("D28") "18.19"
Set-TextCell $ws.Range("E28") "  +0.48%  "

Set-TextCell $ws.Range("D29") "1.933"
Set-TextCell $ws.Range("E29") "  -1.69%  "

Set-TextCell $ws.Range("D30") "1.364"
Set-TextCell $ws.Range("E30") "  -5.01%  "

Set-TextCell $ws.Range("E31") "  -0.20%  "

Set-TextCell $ws.Range("D32") "4.317"
Set-TextCell $ws.Range("E32") "  -1.94%  "

Set-TextCell $ws.Range("D33") "4.064"
Set-TextCell $ws.Range("E33") "  +0.12%  "

Set-TextCell $ws.Range("D34") "0.05229"
Set-TextCell $ws.Range("E34") "  +0.42%  "

Set-TextCell $ws.Range("D35") "1.195"
Set-TextCell $ws.Range("E35") "  +2.28%  "

Set-TextCell $ws.Range("D36") "0.7213"
Set-TextCell $ws.Range("E36") "  +2.00%  "

Set-TextCell $ws.Range("E37") "  -0.11%  "

Set-TextCell $ws.Range("D38") "2.675"
Set-TextCell $ws.Range("E38") "  +0.21%  "

Set-TextCell $ws.Range("E39") "  +0.52%  "

Set-TextCell $ws.Range("D40") "2.697"
Set-TextCell $ws.Range("E40") "  -1.18%  "

Set-TextCell $ws.Range("D41") "1.174.21"
Set-TextCell $ws.Range("E41") "  +3.51%  "

Set-TextCell $ws.Range("D42") "0.9168"
Set-TextCell $ws.Range("E42") "  -0.15%  "

Set-TextCell $ws.Range("D43") "6.001"
Set-TextCell $ws.Range("E43") "  +1.14%  "

Set-TextCell $ws.Range("B44") "TheSandbox"
Set-TextCell $ws.Range("C44") "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws.Range("D44") "0.4307"
Set-TextCell $ws.Range("E44") "  +0.68%  "

Set-TextCell $ws.Range("B45") "Aave"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws.Range("D45") "71.56"
Set-TextCell $ws.Range("E45") "  +1.28%  "

Set-TextCell $ws.Range("D46") "0.9994"
Set-TextCell $ws.Range("E46") "  -0.05%  "

Set-TextCell $ws.Range("D47") "102.37"
Set-TextCell $ws.Range("E47") "  +0.01%  "

Set-TextCell $ws.Range("D48") "0.5311"
Set-TextCell $ws.Range("E48") "  -2.06%  "

Set-TextCell $ws.Range("D49") "1.765"
Set-TextCell $ws.Range("E49") "  -0.28%  "

Set-TextCell $ws.Range("B50") "SynthetixNetwork"
Set-TextCell $ws.Range("C50") "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
Set-TextCell $ws.Range("D50") "2.897"
Set-TextCell $ws.Range("E50") "  +4.74%  "

Set-TextCell $ws.Range("B51") "EnergySwap"
Set-TextCell $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws.Range("D51") "9.205"
Set-TextCell $ws.Range("E51") "  +0.29%  "
